# Applies the commit: adds a new "ODI Bowling Extra" worksheet populated with
# MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL data, mirroring the
# layout/style already used by the "ODI Batting Extra" sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook ----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$ws.Range("A1:C1").Style = "Header"

# --- Data rows -----------------------------------------------------------
$data = @(
    @("4340", "1", "20.00%"),
    @("4348", "0", ""),
    @("4377", "1", "20.00%"),
    @("4378", "1", "10.00%"),
    @("4379", "1", "20.00%"),
    @("4444", "0", "20.00%"),
    @("4446", "0", "30.00%"),
    @("4448", "0", "20.00%"),
    @("4525", "", ""),
    @("4528", "0", "40.00%"),
    @("4530", "0", "10.00%"),
    @("4537", "0", "10.00%"),
    @("4538", "", ""),
    @("4539", "0", ""),
    @("4582", "0", ""),
    @("4585", "0", ""),
    @("4588", "1", "10.00%"),
    @("4671", "", ""),
    @("4674", "0", ""),
    @("4675", "", "")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

$ws.Range("A1").Select()
